$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 193.77730009557
$ws.Range("D2").Value = 0.0000000000000000000000000000000000000000817572480199464
$ws.Range("B3").Value = 4177.64555793754
$ws.Range("B4").Value = 1968.24894106623
